$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Closing_Price (M2SL data) - update trailing observations and
# append a new final row.
# ---------------------------------------------------------------------------
$wsPrice = $wb.Worksheets.Item("Closing_Price")

$priceUpdates = @{
    253 = 19114.3
    254 = 19356.7
    255 = 19600.3
    256 = 19841.2
    257 = 20116.9
    258 = 20431.1
    259 = 20506.6
    260 = 20662.9
    261 = 20847.8
    262 = 20964.3
    263 = 21116.2
    264 = 21316.1
    265 = 21549.1
    266 = 21561.8
    267 = 21570.3
    268 = 21697.5
    269 = 21677.1
    270 = 21665.2
    271 = 21665.7
    272 = 21702.6
    273 = 21658.9
    275 = 21432.3
    276 = 21399
    278 = 21212.7
    279 = 21077.4
    280 = 20841.1
    281 = 20674.6
}

foreach ($rowNum in $priceUpdates.Keys) {
    $wsPrice.Cells.Item($rowNum, 2).Value = $priceUpdates[$rowNum]
}

# Append new row 282 (date 2023-05-01 -> serial 45047, value 20805.5)
$wsPrice.Cells.Item(281, 1).Copy($wsPrice.Cells.Item(282, 1))
$wsPrice.Cells.Item(282, 1).Value = 45047
$wsPrice.Cells.Item(282, 2).Value = 20805.5

# ---------------------------------------------------------------------------
# Sheet 2: SeriesInfo - update metadata fields and append a new "Source" row.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# These B-column values look like dates/timestamps; force the cells to text
# format first so they round-trip as literal strings (matching the
# inlineStr content in the source) instead of being auto-converted to date
# serial numbers.
$wsInfo.Range("B3:B4").NumberFormat = "@"
$wsInfo.Cells.Item(7, 2).NumberFormat = "@"
$wsInfo.Cells.Item(14, 2).NumberFormat = "@"

$wsInfo.Cells.Item(3, 2).Value = "2023-07-13"
$wsInfo.Cells.Item(4, 2).Value = "2023-07-13"
$wsInfo.Cells.Item(7, 2).Value = "2023-05-01"
$wsInfo.Cells.Item(14, 2).Value = "2023-06-27 12:03:02-05"
$wsInfo.Cells.Item(15, 2).Value = 93

$wsInfo.Cells.Item(1, 1).Copy($wsInfo.Cells.Item(17, 1))
$wsInfo.Cells.Item(17, 1).Value = "Source"
$wsInfo.Cells.Item(17, 2).Value = "fred"
